# Actualización desde MV -datos-
# Appends the new daily rows (08-09-2021 .. 15-09-2021) to the bottom of the
# "Spot posiciones netas y suscripciones" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value (e.g. a "dd-mm-yyyy" style date string)
# into a cell without letting Excel auto-convert it into a real date serial
# number (which would also silently create/attach a new number-format style).
# We do this by building the text through a formula (so its result is
# explicitly a string), then copying that computed value with
# PasteSpecial(values-only) onto the destination cell. A value arriving this
# way is stored as plain text, exactly like the original cells in column A.
function Set-LiteralText($address, $text) {
    $helper = $ws.Range("ZZ1")
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $ws.Range($address).PasteSpecial(-4163)
    $helper.Value = ""
    $excel.CutCopyMode = $false
}

$rows = @(
    @{ Row = 174; Date = "08-09-2021"; B = -10806; C = 3020; D = 958; E = 750;  F = 1313 },
    @{ Row = 175; Date = "09-09-2021"; B = -10807; C = 3044; D = 895; E = 861;  F = 1288 },
    @{ Row = 176; Date = "10-09-2021"; B = -10850; C = 2851; D = 825; E = 848;  F = 1178 },
    @{ Row = 177; Date = "13-09-2021"; B = -11127; C = 2788; D = 748; E = 1003; F = 1037 },
    @{ Row = 178; Date = "14-09-2021"; B = -11105; C = 3796; D = 978; E = 941;  F = 1878 },
    @{ Row = 179; Date = "15-09-2021"; B = -11103; C = 4092; D = 964; E = 957;  F = 2170 }
)

foreach ($r in $rows) {
    Set-LiteralText ("A" + $r.Row) $r.Date
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
}
